# ComPADRE_tracker_experiments.xlsx — mark completed tests
#
# Commit message: "updated with completed tests"
#
# Columns: A=Online?, B=title, C=Sync OK?, D=mp4 or image?, E=ver5 OK?,
#          F=Java Notes, G=JavaScript Notes
#
# For a large batch of rows, the "mp4 or image?" (D) and "ver5 OK?" (E)
# tests are now marked complete ("yes"), and row 8 — which was still
# highlighted/custom-formatted as "pending" — has that formatting cleared
# now that it, too, is complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsWithDandE = @(4,5,7,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25,26,27,29,30,31,32,33,34,35,36,37,38)
foreach ($r in $rowsWithDandE) {
    $ws.Range("D$r").Value = "yes"
    $ws.Range("E$r").Value = "yes"
}

$rowsWithDOnly = @(48)
foreach ($r in $rowsWithDOnly) {
    $ws.Range("D$r").Value = "yes"
}

# Row 8 was the last remaining "still pending" row and carried a custom
# highlight format (style applied to A8/C8/G8 plus the row itself); clear
# it now that its tests are complete, but keep B8's hyperlink style intact.
$ws.Rows.Item(8).ClearFormats()
$ws.Range("B8").Style = "Hyperlink"

# Reflect the reviewer's new focus area in the UI state: frozen pane/
# scroll position and selection moved down to the newly completed block.
$ws.Range("A11").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("D3:E10").Select()
